$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: bold the lowercase "responsive" in
#   "...sida på olika sätt. I responsive web design har det fokuserats..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("I responsive web design har", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordRng = $d.Range($rng.Start, $rng.End)
$wordRng.Find.Execute("responsive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordRng.Bold = 1

# ---------------------------------------------------------------------------
# Edit 2: "Iphone" -> "iPhone", and bold "smartphones" in
#   "...utvecklingen för smartphones eskalerat markant."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Iphone", $true, $false, $false, $false, $false, $true, 1, $false, "iPhone", 2)

$rng2 = $d.Content
$rng2.Find.Execute("har utvecklingen för smartphones eskalerat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$smartRng = $d.Range($rng2.Start, $rng2.End)
$smartRng.Find.Execute("smartphones", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$smartRng.Bold = 1

# ---------------------------------------------------------------------------
# Edit 3: " använder mobilt internet" -> " använder sig av mobilt internet"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("använder mobilt internet", $true, $false, $false, $false, $false, $true, 1, $false, "använder sig av mobilt internet", 2)

# ---------------------------------------------------------------------------
# Edit 4: rewrite the final two body paragraphs
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("helt och hållet. Datorer har fortfarande", $true, $false, $false, $false, $false, $true, 1, $false, "helt och hållet, vilket till en början var tanken. Datorer har fortfarande", 2)

$d.Content.Find.Execute("Även om mobila användarna enligt statistik kommer vara flera än desktop användarna så går det inte att ignorera dessa användare då dom flesta är både och. ", $true, $false, $false, $false, $false, $true, 1, $false, "Även om mobila användarna enligt statistik kommer vara flera än desktop användarna så går det inte att ignorera dessa användare då de flesta är användare utav båda.", 2)
